$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.543.01'
$ws.Range("E2").Value = '  -4.47%  '
$ws.Range("D3").Value = '3.333.68'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.12'
$ws.Range("E5").Value = '  -3.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.39'
$ws.Range("E6").Value = '  -5.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +3.27%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -3.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("E11").Value = '  -4.14%  '
$ws.Range("D12").Value = '3.913.30'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.99'
$ws.Range("E14").Value = '  -6.24%  '
$ws.Range("D15").Value = '66.662.00'
$ws.Range("E15").Value = '  -4.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000167'
$ws.Range("E16").Value = '  -2.73%  '
$ws.Range("D17").Value = '3.352.74'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '436.78'
$ws.Range("E18").Value = '  -3.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.68'
$ws.Range("E19").Value = '  -2.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.55'
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.60'
$ws.Range("E21").Value = '  -3.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.48'
$ws.Range("E22").Value = '  -3.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.517'
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000117'
$ws.Range("E25").Value = '  -4.42%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.02'
$ws.Range("E27").Value = '  -5.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  -2.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.81'
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.28'
$ws.Range("E32").Value = '  -6.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.76'
$ws.Range("E33").Value = '  -3.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.22'
$ws.Range("E34").Value = '  -5.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '162.73'
$ws.Range("E35").Value = '  -1.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.48'
$ws.Range("E36").Value = '  -6.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.41'
$ws.Range("E37").Value = '  -3.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.83'
$ws.Range("E38").Value = '  -6.02%  '
$ws.Range("D39").Value = '2.813.77'
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.793'
$ws.Range("E40").Value = '  -3.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.42'
$ws.Range("E41").Value = '  -4.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.20'
$ws.Range("E42").Value = '  -6.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.17'
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0666'
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.41'
$ws.Range("E45").Value = '  -4.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.35'
$ws.Range("E46").Value = '  -7.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '320.43'
$ws.Range("E47").Value = '  -6.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0273'
$ws.Range("E48").Value = '  -4.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.103'
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.977'
$ws.Range("E50").Value = '  -4.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.15'
$ws.Range("E51").Value = '  -3.05%  '
